# Auto-generated edit script: update '想去人数' (want-to-go count) figures
# per the commit's refreshed bilibili scrape snapshot.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F3").Value = 27012
$ws1.Range("G4").Value = 78
$ws1.Range("F5").Value = 641
$ws1.Range("F6").Value = 189
$ws1.Range("F7").Value = 562
$ws1.Range("F9").Value = 374
$ws1.Range("F10").Value = 473
$ws1.Range("F13").Value = 313
$ws1.Range("F14").Value = 95
$ws1.Range("F15").Value = 475
$ws1.Range("F17").Value = 1613
$ws1.Range("F18").Value = 245
$ws1.Range("F19").Value = 493
$ws1.Range("F20").Value = 136
$ws1.Range("F21").Value = 454

# Sheet 2
$ws2.Range("F2").Value = 4521
$ws2.Range("F3").Value = 245
$ws2.Range("F9").Value = 3
$ws2.Range("F10").Value = 116
$ws2.Range("F11").Value = 450
$ws2.Range("F17").Value = 74
$ws2.Range("F18").Value = 27
$ws2.Range("F20").Value = 29

# Sheet 3
$ws3.Range("F2").Value = 5170

# Sheet 4
$ws4.Range("F3").Value = 5170
$ws4.Range("F5").Value = 27012
$ws4.Range("F6").Value = 4521
$ws4.Range("G7").Value = 78
$ws4.Range("F8").Value = 245
$ws4.Range("F9").Value = 641
$ws4.Range("F12").Value = 189
$ws4.Range("F16").Value = 3
$ws4.Range("F17").Value = 116
$ws4.Range("F18").Value = 450
$ws4.Range("F19").Value = 562
$ws4.Range("F23").Value = 374
$ws4.Range("F24").Value = 473
$ws4.Range("F28").Value = 313
$ws4.Range("F29").Value = 95
$ws4.Range("F32").Value = 475
$ws4.Range("F34").Value = 74
$ws4.Range("F35").Value = 1613
$ws4.Range("F36").Value = 245
$ws4.Range("F37").Value = 494
$ws4.Range("F38").Value = 27
$ws4.Range("F39").Value = 136
$ws4.Range("F40").Value = 454
$ws4.Range("F45").Value = 29
